$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "A low quality drawing feature ..." -> "A low-quality drawing feature ..."
#    Scope the Find/Replace to the specific paragraph so the other,
#    unrelated occurrence of "low quality" later in the document is
#    left untouched.
# ------------------------------------------------------------------
$para = $d.Paragraphs(6).Range
$para.Find.Execute("low quality", $true, $false, $false, $false, $false, $true, 1, $false, "low-quality", 2)

# ------------------------------------------------------------------
# 2. The text-replace above normalizes (coalesces) every run inside the
#    paragraph into a single run. Re-create the expected run
#    boundaries ("A " | "low-quality" | " drawing feature...set up. " |
#    "An offline database...") by toggling a character property across
#    exactly the desired sub-ranges; toggling it back to its original
#    value keeps the formatting unchanged while forcing Word to split
#    the run at those offsets.
# ------------------------------------------------------------------
$full = $d.Content.Text
$idxA = $full.IndexOf("A low-quality")
$idxLowQuality = $idxA + 2
$idxAfterLowQuality = $idxLowQuality + ("low-quality").Length
$idxAnOffline = $full.IndexOf("An offline database")
$idxNextPara = $full.IndexOf("All the screens")

$rWord = $d.Range($idxLowQuality, $idxAfterLowQuality)
$rWord.Bold = 1
$rWord.Bold = 0

$rTail = $d.Range($idxAnOffline, $idxNextPara)
$rTail.Bold = 1
$rTail.Bold = 0

# ------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the document (last
#    edit location before this change) to sit right after "A ", i.e.
#    immediately before "low-quality" - this mirrors Word recording
#    the new last-edit position.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($idxLowQuality, $idxLowQuality))

Write-Output "done"
